# Add "Start" and "Develop" buttons to the Translation sheet of the
# texts.xlsx workbook. This adds two new rows (7 and 8) to the
# "Translation" worksheet, each describing a new UI text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 7: SingleUseId3 / Default / Center / LTR / Start
$ws.Range("B7").Value = "SingleUseId3"
$ws.Range("C7").Value = "Default"
$ws.Range("D7").Value = "Center"
$ws.Range("E7").Value = "LTR"
$ws.Range("F7").Value = "Start"

# Row 8: SingleUseId4 / Default / Center / LTR / Develop
$ws.Range("B8").Value = "SingleUseId4"
$ws.Range("C8").Value = "Default"
$ws.Range("D8").Value = "Center"
$ws.Range("E8").Value = "LTR"
$ws.Range("F8").Value = "Develop"
